# Update countries & provincias Spain
#
# Refreshes the COVID snapshot on sheet "Pais":
#   - bumps the "Datos actualizados ..." timestamp in A1
#   - updates Casos totales/Nuevos casos/Casos activos/Recuperados/
#     Casos criticos/Muertes hoy/Muertes (columns B-H) for several
#     countries whose figures moved
#   - a handful of neighbouring countries swapped rank order, which
#     (since country names are looked up by row position) shows up as
#     the country label in column A changing for that row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 12:05"

# Row 4 - Estados Unidos (updated figures only)
$ws.Range("B4").Value = 1571018
$ws.Range("C4").Value = 435
$ws.Range("D4").Value = 361227
$ws.Range("E4").Value = 1116249
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 93542

# Row 37 - was Sudafrica, now Rumania (rank swap + updated figures)
$ws.Range("A37").Value = "Rumania"
$ws.Range("B37").Value = 17387
$ws.Range("C37").Value = 196
$ws.Range("D37").Value = 10356
$ws.Range("E37").Value = 5890
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 1141

# Row 38 - was Rumania, now Sudafrica (rank swap)
$ws.Range("A38").Value = "Sudafrica"
$ws.Range("B38").Value = 17200
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 7960
$ws.Range("E38").Value = 8928
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 312

# Row 65 - Oman (updated figures only)
$ws.Range("B65").Value = 6043
$ws.Range("C65").Value = 372
$ws.Range("D65").Value = 1661
$ws.Range("E65").Value = 4355
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 27

# Row 89 - Estonia (updated figures only)
$ws.Range("D89").Value = 956
$ws.Range("E89").Value = 774

# Row 102 - Hong Kong (updated figures only)
$ws.Range("D102").Value = 1026
$ws.Range("E102").Value = 26

# Row 107 - was Kenia, now Albania (rank swap + updated figures)
$ws.Range("A107").Value = "Albania"
$ws.Range("B107").Value = 964
$ws.Range("C107").Value = 15
$ws.Range("D107").Value = 758
$ws.Range("E107").Value = 175
$ws.Range("H107").Value = 31

# Row 108 - was Libano, now Kenia (rank swap)
$ws.Range("A108").Value = "Kenia"
$ws.Range("B108").Value = 963
$ws.Range("D108").Value = 358
$ws.Range("E108").Value = 555
$ws.Range("H108").Value = 50

# Row 109 - was Albania, now Libano (rank swap)
$ws.Range("A109").Value = "Libano"
$ws.Range("B109").Value = 954
$ws.Range("D109").Value = 251
$ws.Range("E109").Value = 677
$ws.Range("H109").Value = 26

# Row 136 - was Republica de Africa Central, now Etiopia (rank swap + updated figures)
$ws.Range("A136").Value = "Etiopia"
$ws.Range("B136").Value = 389
$ws.Range("C136").Value = 24
$ws.Range("D136").Value = 122
$ws.Range("E136").Value = 262
$ws.Range("H136").Value = 5

# Row 137 - was Etiopia, now Republica de Africa Central (rank swap)
$ws.Range("A137").Value = "Republica de Africa Central"
$ws.Range("B137").Value = 366
$ws.Range("D137").Value = 18
$ws.Range("E137").Value = 348
$ws.Range("H137").Value = 0

# Row 196 - was Nueva Caledonia, now Santa Lucia (rank swap, figures tie)
$ws.Range("A196").Value = "Santa Lucia"

# Row 197 - was Santa Lucia, now Nueva Caledonia (rank swap, figures tie)
$ws.Range("A197").Value = "Nueva Caledonia"

# Row 209 - was Montserrat, now Groenlandia (rank swap + updated figures)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

# Row 210 - was Groenlandia, now Seychelles (rank swap, figures tie)
$ws.Range("A210").Value = "Seychelles"

# Row 211 - was Seychelles, now Montserrat (rank swap + updated figures)
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Row 214 - was San Bartolome, now Bonaire, San Eustaquio y Saba (rank swap, figures tie)
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"

# Row 216 - was Bonaire, San Eustaquio y Saba, now San Bartolome (rank swap, figures tie)
$ws.Range("A216").Value = "San Bartolome"
